$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("B2").Value = 1.02
$ws.Range("C2").Value = 1.053311405987383
$ws.Range("D2").Value = 1.057178970700206
$ws.Range("E2").Value = 1.04994373469313
$ws.Range("F2").Value = 1.065863333543118
$ws.Range("I2").Value = 1.042071305733029
$ws.Range("J2").Value = 1.058329415949934
$ws.Range("K2").Value = 1.059914650324083
$ws.Range("L2").Value = 1.052699386186681
$ws.Range("M2").Value = 1.068575428680131
$ws.Range("B3").Value = 1.02
$ws.Range("C3").Value = 1.054917820275471
$ws.Range("D3").Value = 1.058435077303894
$ws.Range("E3").Value = 1.051334340832682
$ws.Range("F3").Value = 1.067255699464447
$ws.Range("I3").Value = 1.042456115581319
$ws.Range("J3").Value = 1.05958347505465
$ws.Range("K3").Value = 1.060983418105704
$ws.Range("L3").Value = 1.053900848627291
$ws.Range("M3").Value = 1.06978183564112
$ws.Range("B4").Value = 1.02
$ws.Range("C4").Value = 1.055955591618639
$ws.Range("D4").Value = 1.05924616284633
$ws.Range("E4").Value = 1.052232808980301
$ws.Range("F4").Value = 1.068155191664532
$ws.Range("I4").Value = 1.042702974607754
$ws.Range("J4").Value = 1.060392890346863
$ws.Range("K4").Value = 1.061672727730932
$ws.Range("L4").Value = 1.054676404003347
$ws.Range("M4").Value = 1.070560469524717
$ws.Range("B5").Value = 1.02
$ws.Range("C5").Value = 1.05639147645416
$ws.Range("D5").Value = 1.059586742742082
$ws.Range("E5").Value = 1.052610209725866
$ws.Range("F5").Value = 1.068532995243461
$ws.Range("I5").Value = 1.042806244618098
$ws.Range("J5").Value = 1.060732686084891
$ws.Range("K5").Value = 1.061961979698917
$ws.Range("L5").Value = 1.05500200615906
$ws.Range("M5").Value = 1.070887336739887
$ws.Range("B6").Value = 1.02
$ws.Range("C6").Value = 1.056464640590529
$ws.Range("D6").Value = 1.059643904306206
$ws.Range("E6").Value = 1.052673558693224
$ws.Range("F6").Value = 1.068596410218494
$ws.Range("I6").Value = 1.042823554300172
$ws.Range("J6").Value = 1.060789711198179
$ws.Range("K6").Value = 1.062010515247211
$ws.Range("L6").Value = 1.05505665050625
$ws.Range("M6").Value = 1.070942191781378
$ws.Range("B7").Value = 1.02
$ws.Range("C7").Value = 1.05596141747059
$ws.Range("D7").Value = 1.059250715257258
$ws.Range("E7").Value = 1.0522378530574
$ws.Range("F7").Value = 1.068160241233356
$ws.Range("I7").Value = 1.042704356505721
$ws.Range("J7").Value = 1.060397432601119
$ws.Range("K7").Value = 1.061676594820013
$ws.Range("L7").Value = 1.05468075644121
$ws.Range("M7").Value = 1.070564838981128
$ws.Range("B8").Value = 1.02
$ws.Range("C8").Value = 1.053854656027642
$ws.Range("D8").Value = 1.057603833537508
$ws.Range("E8").Value = 1.050413979811986
$ws.Range("F8").Value = 1.066334196986565
$ws.Range("I8").Value = 1.04220179834801
$ws.Range("J8").Value = 1.058753658761089
$ws.Range("K8").Value = 1.060276315826143
$ws.Range("L8").Value = 1.053105817487065
$ws.Range("M8").Value = 1.068983555926171
$ws.Range("B9").Value = 1.02
$ws.Range("C9").Value = 1.050128935912557
$ws.Range("D9").Value = 1.054688514744175
$ws.Range("E9").Value = 1.04718944107112
$ws.Range("F9").Value = 1.063104964945842
$ws.Range("I9").Value = 1.041299735161952
$ws.Range("J9").Value = 1.055841139911728
$ws.Range("K9").Value = 1.057791309623664
$ws.Range("L9").Value = 1.050315944357748
$ws.Range("M9").Value = 1.066181586151517
$ws.Range("B10").Value = 1.02
$ws.Range("C10").Value = 1.047635560798287
$ws.Range("D10").Value = 1.052735620548487
$ws.Range("E10").Value = 1.04503213792172
$ws.Range("F10").Value = 1.06094396368233
$ws.Range("I10").Value = 1.040687113983104
$ws.Range("J10").Value = 1.053888269292959
$ws.Range("K10").Value = 1.056122471364848
$ws.Range("L10").Value = 1.048445773111025
$ws.Range("M10").Value = 1.064302745382207
$ws.Range("B11").Value = 1.02
$ws.Range("C11").Value = 1.046553501044236
$ws.Range("D11").Value = 1.051887687750231
$ws.Range("E11").Value = 1.044096094475273
$ws.Range("F11").Value = 1.060006186959753
$ws.Range("I11").Value = 1.040419139359406
$ws.Range("J11").Value = 1.053039899970462
$ws.Range("K11").Value = 1.055396875304433
$ws.Range("L11").Value = 1.04763344168309
$ws.Range("M11").Value = 1.063486518978061
$ws.Range("B12").Value = 1.02
$ws.Range("C12").Value = 1.046151201990267
$ws.Range("D12").Value = 1.051572371780543
$ws.Range("E12").Value = 1.043748109373132
$ws.Range("F12").Value = 1.059657538469735
$ws.Range("I12").Value = 1.040319192035509
$ws.Range("J12").Value = 1.052724353919394
$ws.Range("K12").Value = 1.055126901805717
$ws.Range("L12").Value = 1.047331316585568
$ws.Range("M12").Value = 1.063182926136689
$ws.Range("B13").Value = 1.02
$ws.Range("C13").Value = 1.046237513608037
$ws.Range("D13").Value = 1.051640024362796
$ws.Range("E13").Value = 1.043822766937197
$ws.Range("F13").Value = 1.059732339214326
$ws.Range("I13").Value = 1.04034064966236
$ws.Range("J13").Value = 1.052792058942385
$ws.Range("K13").Value = 1.055184832738346
$ws.Range("L13").Value = 1.047396141181545
$ws.Range("M13").Value = 1.063248066514555
$ws.Range("B14").Value = 1.02
$ws.Range("C14").Value = 1.046520254573048
$ws.Range("D14").Value = 1.051861630941608
$ws.Range("E14").Value = 1.044067336026317
$ws.Range("F14").Value = 1.059977374072733
$ws.Range("I14").Value = 1.040410886056708
$ws.Range("J14").Value = 1.053013825533894
$ws.Range("K14").Value = 1.055374568531891
$ws.Range("L14").Value = 1.047608475912137
$ws.Range("M14").Value = 1.063461432309638
$ws.Range("B15").Value = 1.02
$ws.Range("C15").Value = 1.046694410723696
$ws.Range("D15").Value = 1.05199812266078
$ws.Range("E15").Value = 1.044217983506022
$ws.Range("F15").Value = 1.060128305942428
$ws.Range("I15").Value = 1.040454106644991
$ws.Range("J15").Value = 1.053150406815261
$ws.Range("K15").Value = 1.055491410526736
$ws.Range("L15").Value = 1.04773925055027
$ws.Range("M15").Value = 1.063592839447355
$ws.Range("B16").Value = 1.02
$ws.Range("C16").Value = 1.047707321743474
$ws.Range("D16").Value = 1.052791845571244
$ws.Range("E16").Value = 1.045094218897633
$ws.Range("F16").Value = 1.061006156931386
$ws.Range("I16").Value = 1.040704841305751
$ws.Range("J16").Value = 1.053944513804329
$ws.Range("K16").Value = 1.056170563411529
$ws.Range("L16").Value = 1.048499630766869
$ws.Range("M16").Value = 1.064356858581653
$ws.Range("B17").Value = 1.02
$ws.Range("C17").Value = 1.048342040231026
$ws.Range("D17").Value = 1.053289101166261
$ws.Range("E17").Value = 1.045643338736224
$ws.Range("F17").Value = 1.061556255039519
$ws.Range("I17").Value = 1.040861393964888
$ws.Range("J17").Value = 1.05444189014093
$ws.Range("K17").Value = 1.056595775558403
$ws.Range("L17").Value = 1.048975912752993
$ws.Range("M17").Value = 1.064835385558804
$ws.Range("B18").Value = 1.02
$ws.Range("C18").Value = 1.04871202885899
$ws.Range("D18").Value = 1.053578918950358
$ws.Range("E18").Value = 1.045963446695054
$ws.Range("F18").Value = 1.061876920856236
$ws.Range("I18").Value = 1.040952447643325
$ws.Range("J18").Value = 1.054731735303781
$ws.Range("K18").Value = 1.056843507771146
$ws.Range("L18").Value = 1.049253475950168
$ws.Range("M18").Value = 1.065114244601194
$ws.Range("B19").Value = 1.02
$ws.Range("C19").Value = 1.048838146427919
$ws.Range("D19").Value = 1.053677701739976
$ws.Range("E19").Value = 1.046072564234804
$ws.Range("F19").Value = 1.061986226440728
$ws.Range("I19").Value = 1.04098345044593
$ws.Range("J19").Value = 1.054830520151985
$ws.Range("K19").Value = 1.056927929623324
$ws.Range("L19").Value = 1.049348076682767
$ws.Range("M19").Value = 1.065209284866954
$ws.Range("B20").Value = 1.02
$ws.Range("C20").Value = 1.048273965016969
$ws.Range("D20").Value = 1.053235773469267
$ws.Range("E20").Value = 1.045584442507192
$ws.Range("F20").Value = 1.061497255168399
$ws.Range("I20").Value = 1.040844624348633
$ws.Range("J20").Value = 1.054388553948059
$ws.Range("K20").Value = 1.056550184041789
$ws.Range("L20").Value = 1.048924837503561
$ws.Range("M20").Value = 1.064784070850359
$ws.Range("B21").Value = 1.02
$ws.Range("C21").Value = 1.046437004820017
$ws.Range("D21").Value = 1.051796383218079
$ws.Range("E21").Value = 1.043995324792653
$ws.Range("F21").Value = 1.059905226219854
$ws.Range("I21").Value = 1.040390214521884
$ws.Range("J21").Value = 1.052948532572612
$ws.Range("K21").Value = 1.055318708684535
$ws.Range("L21").Value = 1.047545959431715
$ws.Range("M21").Value = 1.063398612783784
$ws.Range("B22").Value = 1.02
$ws.Range("C22").Value = 1.0452798660869
$ws.Range("D22").Value = 1.050889317909541
$ws.Range("E22").Value = 1.042994461247335
$ws.Range("F22").Value = 1.058902419155564
$ws.Range("I22").Value = 1.040102137538738
$ws.Range("J22").Value = 1.052040676980054
$ws.Range("K22").Value = 1.054541796123095
$ws.Range("L22").Value = 1.046676749108441
$ws.Range("M22").Value = 1.062525144098507
$ws.Range("B23").Value = 1.02
$ws.Range("C23").Value = 1.045893496686763
$ws.Range("D23").Value = 1.051370368650212
$ws.Range("E23").Value = 1.043525204178384
$ws.Range("F23").Value = 1.059434203019125
$ws.Range("I23").Value = 1.04025507844012
$ws.Range("J23").Value = 1.052522184236784
$ws.Range("K23").Value = 1.0549539044357
$ws.Range("L23").Value = 1.047137750442718
$ws.Range("M23").Value = 1.062988414294826
$ws.Range("B24").Value = 1.02
$ws.Range("C24").Value = 1.048304725992697
$ws.Range("D24").Value = 1.053259870649298
$ws.Range("E24").Value = 1.045611055747751
$ws.Range("F24").Value = 1.061523915279779
$ws.Range("I24").Value = 1.040852202622431
$ws.Range("J24").Value = 1.054412655099977
$ws.Range("K24").Value = 1.056570785772145
$ws.Range("L24").Value = 1.048947916964448
$ws.Range("M24").Value = 1.064807258554575
$ws.Range("B25").Value = 1.02
$ws.Range("C25").Value = 1.051093765322707
$ws.Range("D25").Value = 1.055443812387914
$ws.Range("E25").Value = 1.048024369385662
$ws.Range("F25").Value = 1.063941207962688
$ws.Range("I25").Value = 1.041534910598016
$ws.Range("J25").Value = 1.056596035352925
$ws.Range("K25").Value = 1.058435860021317
$ws.Range("L25").Value = 1.051038970107731
$ws.Range("M25").Value = 1.066907846690206
